$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.3
$ws.Cells.Item(2, 9).Value = 3.4
$ws.Cells.Item(2, 10).Value = 1.08
$ws.Cells.Item(2, 12).Value = 1.36
$ws.Cells.Item(2, 14).Value = 2.2
$ws.Cells.Item(2, 15).Value = 1.65
$ws.Cells.Item(2, 20).Value = 7
$ws.Cells.Item(2, 23).Value = 21

# Row 3
$ws.Cells.Item(3, 7).Value = 1.75
$ws.Cells.Item(3, 8).Value = 3.3
$ws.Cells.Item(3, 9).Value = 5.25
$ws.Cells.Item(3, 10).Value = 1.11
$ws.Cells.Item(3, 11).Value = 6.5
$ws.Cells.Item(3, 12).Value = 1.53
$ws.Cells.Item(3, 13).Value = 2.38
$ws.Cells.Item(3, 14).Value = 2.7
$ws.Cells.Item(3, 15).Value = 1.44

# Row 4
$ws.Cells.Item(4, 7).Value = 2.32
$ws.Cells.Item(4, 10).Value = 1.1
$ws.Cells.Item(4, 12).Value = 1.44
$ws.Cells.Item(4, 13).Value = 2.63
$ws.Cells.Item(4, 16).Value = 1.54

# Row 5
$ws.Cells.Item(5, 7).Value = 2.2
$ws.Cells.Item(5, 8).Value = 3.1
$ws.Cells.Item(5, 9).Value = 3.4
$ws.Cells.Item(5, 10).Value = 1.08
$ws.Cells.Item(5, 11).Value = 8
$ws.Cells.Item(5, 12).Value = 1.4
$ws.Cells.Item(5, 13).Value = 2.75
$ws.Cells.Item(5, 14).Value = 2.35
$ws.Cells.Item(5, 15).Value = 1.57
$ws.Cells.Item(5, 16).Value = 1.5
$ws.Cells.Item(5, 17).Value = 2.37
$ws.Cells.Item(5, 20).Value = 6.5
$ws.Cells.Item(5, 21).Value = 10
$ws.Cells.Item(5, 22).Value = 9.5
$ws.Cells.Item(5, 23).Value = 21
$ws.Cells.Item(5, 24).Value = 21
$ws.Cells.Item(5, 26).Value = 7.5
$ws.Cells.Item(5, 27).Value = 6
$ws.Cells.Item(5, 28).Value = 17
$ws.Cells.Item(5, 29).Value = 51
$ws.Cells.Item(5, 31).Value = 8.5
$ws.Cells.Item(5, 32).Value = 15
$ws.Cells.Item(5, 33).Value = 13
$ws.Cells.Item(5, 34).Value = 34
$ws.Cells.Item(5, 35).Value = 29
$ws.Cells.Item(5, 36).Value = 41

# Row 6
$ws.Cells.Item(6, 7).Value = 1.66
$ws.Cells.Item(6, 8).Value = 3.7
$ws.Cells.Item(6, 9).Value = 4.75
$ws.Cells.Item(6, 10).Value = 1.05
$ws.Cells.Item(6, 12).Value = 1.29
$ws.Cells.Item(6, 13).Value = 3.5
$ws.Cells.Item(6, 14).Value = 1.93
$ws.Cells.Item(6, 15).Value = 1.93
$ws.Cells.Item(6, 16).Value = 1.33
$ws.Cells.Item(6, 31).Value = 13
$ws.Cells.Item(6, 33).Value = 15

# Row 7
$ws.Cells.Item(7, 9).Value = 2.82
$ws.Cells.Item(7, 10).Value = 1.1
$ws.Cells.Item(7, 11).Value = 7
$ws.Cells.Item(7, 12).Value = 1.5
$ws.Cells.Item(7, 16).Value = 1.54

# Row 8
$ws.Cells.Item(8, 7).Value = 1.62
$ws.Cells.Item(8, 16).Value = 1.41
$ws.Cells.Item(8, 17).Value = 2.62

# Row 9
$ws.Cells.Item(9, 7).Value = 1.38
$ws.Cells.Item(9, 8).Value = 4.2
$ws.Cells.Item(9, 9).Value = 8.5
$ws.Cells.Item(9, 10).Value = 1.05
$ws.Cells.Item(9, 11).Value = 11
$ws.Cells.Item(9, 12).Value = 1.3
$ws.Cells.Item(9, 13).Value = 3.4
$ws.Cells.Item(9, 14).Value = 2
$ws.Cells.Item(9, 15).Value = 1.85
$ws.Cells.Item(9, 16).Value = 1.4
$ws.Cells.Item(9, 17).Value = 2.75
$ws.Cells.Item(9, 18).Value = 2.25
$ws.Cells.Item(9, 19).Value = 1.57
$ws.Cells.Item(9, 20).Value = 6
$ws.Cells.Item(9, 21).Value = 6
$ws.Cells.Item(9, 24).Value = 13
$ws.Cells.Item(9, 25).Value = 34
$ws.Cells.Item(9, 26).Value = 9
$ws.Cells.Item(9, 27).Value = 8.5
$ws.Cells.Item(9, 28).Value = 23
$ws.Cells.Item(9, 29).Value = 81
$ws.Cells.Item(9, 31).Value = 17
$ws.Cells.Item(9, 33).Value = 26
$ws.Cells.Item(9, 35).Value = 67
$ws.Cells.Item(9, 36).Value = 67

# Row 10
$ws.Cells.Item(10, 7).Value = 1.4
$ws.Cells.Item(10, 8).Value = 4.1
$ws.Cells.Item(10, 9).Value = 9.5
$ws.Cells.Item(10, 10).Value = 1.08
$ws.Cells.Item(10, 11).Value = 8
$ws.Cells.Item(10, 14).Value = 2.25
$ws.Cells.Item(10, 15).Value = 1.62
$ws.Cells.Item(10, 23).Value = 8.5
$ws.Cells.Item(10, 26).Value = 7.5
$ws.Cells.Item(10, 27).Value = 8.5
$ws.Cells.Item(10, 28).Value = 29
$ws.Cells.Item(10, 31).Value = 17
$ws.Cells.Item(10, 33).Value = 29
$ws.Cells.Item(10, 34).Value = 126

# Row 11
$ws.Cells.Item(11, 7).Value = 1.9
$ws.Cells.Item(11, 8).Value = 3.3
$ws.Cells.Item(11, 11).Value = 9
$ws.Cells.Item(11, 12).Value = 1.4
$ws.Cells.Item(11, 13).Value = 2.75
$ws.Cells.Item(11, 14).Value = 2.25
$ws.Cells.Item(11, 15).Value = 1.62
$ws.Cells.Item(11, 16).Value = 1.5
$ws.Cells.Item(11, 17).Value = 2.5
$ws.Cells.Item(11, 18).Value = 2.1
$ws.Cells.Item(11, 19).Value = 1.67
$ws.Cells.Item(11, 20).Value = 6
$ws.Cells.Item(11, 22).Value = 9
$ws.Cells.Item(11, 25).Value = 34
$ws.Cells.Item(11, 26).Value = 7.5
$ws.Cells.Item(11, 28).Value = 19
$ws.Cells.Item(11, 29).Value = 67
$ws.Cells.Item(11, 31).Value = 9.5
$ws.Cells.Item(11, 32).Value = 21
$ws.Cells.Item(11, 36).Value = 41

# Row 12
$ws.Cells.Item(12, 8).Value = 3.6
$ws.Cells.Item(12, 9).Value = 6.1
$ws.Cells.Item(12, 12).Value = 1.37
$ws.Cells.Item(12, 13).Value = 2.62
$ws.Cells.Item(12, 14).Value = 2.07
$ws.Cells.Item(12, 15).Value = 1.6
$ws.Cells.Item(12, 18).Value = 2.12
$ws.Cells.Item(12, 19).Value = 1.57
$ws.Cells.Item(12, 20).Value = 5.3
$ws.Cells.Item(12, 21).Value = 6.2
$ws.Cells.Item(12, 22).Value = 8.5
$ws.Cells.Item(12, 23).Value = 10.5
$ws.Cells.Item(12, 24).Value = 14
$ws.Cells.Item(12, 25).Value = 37
$ws.Cells.Item(12, 26).Value = 7.8
$ws.Cells.Item(12, 27).Value = 7.3
$ws.Cells.Item(12, 28).Value = 22
$ws.Cells.Item(12, 29).Value = 150
$ws.Cells.Item(12, 31).Value = 12.5
$ws.Cells.Item(12, 32).Value = 35
$ws.Cells.Item(12, 35).Value = 90

# Row 13
$ws.Cells.Item(13, 7).Value = 2.1
$ws.Cells.Item(13, 9).Value = 3.5
$ws.Cells.Item(13, 11).Value = 9.5
$ws.Cells.Item(13, 20).Value = 7.5
$ws.Cells.Item(13, 21).Value = 10
$ws.Cells.Item(13, 23).Value = 19
$ws.Cells.Item(13, 25).Value = 26
$ws.Cells.Item(13, 30).Value = 201
$ws.Cells.Item(13, 31).Value = 11
$ws.Cells.Item(13, 33).Value = 13
$ws.Cells.Item(13, 34).Value = 41
$ws.Cells.Item(13, 35).Value = 29

# Row 15
$ws.Cells.Item(15, 8).Value = 3.9
$ws.Cells.Item(15, 15).Value = 2.32
$ws.Cells.Item(15, 20).Value = 16.5
$ws.Cells.Item(15, 21).Value = 28
$ws.Cells.Item(15, 23).Value = 65
$ws.Cells.Item(15, 24).Value = 35
$ws.Cells.Item(15, 25).Value = 32
$ws.Cells.Item(15, 32).Value = 10.5
$ws.Cells.Item(15, 36).Value = 17.5

# Row 19
$ws.Cells.Item(19, 7).Value = 1.91
$ws.Cells.Item(19, 9).Value = 3.6
$ws.Cells.Item(19, 10).Value = 1.02
$ws.Cells.Item(19, 11).Value = 12
$ws.Cells.Item(19, 12).Value = 1.25
$ws.Cells.Item(19, 13).Value = 3.75
$ws.Cells.Item(19, 16).Value = 1.36
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = 1.73
$ws.Cells.Item(19, 19).Value = 2
$ws.Cells.Item(19, 20).Value = 8.5
$ws.Cells.Item(19, 31).Value = 12
$ws.Cells.Item(19, 35).Value = 29

# Row 20
$ws.Cells.Item(20, 9).Value = 11
$ws.Cells.Item(20, 16).Value = 1.14
$ws.Cells.Item(20, 17).Value = 5
$ws.Cells.Item(20, 18).Value = 1.67
$ws.Cells.Item(20, 19).Value = 2.1
$ws.Cells.Item(20, 21).Value = 10
$ws.Cells.Item(20, 23).Value = 9
$ws.Cells.Item(20, 24).Value = 10
$ws.Cells.Item(20, 25).Value = 21
$ws.Cells.Item(20, 26).Value = 34
$ws.Cells.Item(20, 28).Value = 21
$ws.Cells.Item(20, 29).Value = 41
$ws.Cells.Item(20, 32).Value = 51
$ws.Cells.Item(20, 34).Value = 126
$ws.Cells.Item(20, 35).Value = 51

# Row 21
$ws.Cells.Item(21, 7).Value = 1.57
$ws.Cells.Item(21, 8).Value = 4.1
$ws.Cells.Item(21, 9).Value = 5
$ws.Cells.Item(21, 10).Value = 1.03
$ws.Cells.Item(21, 11).Value = 10.5
$ws.Cells.Item(21, 12).Value = 1.14
$ws.Cells.Item(21, 14).Value = 1.53
$ws.Cells.Item(21, 15).Value = 2.38
$ws.Cells.Item(21, 16).Value = 1.29
$ws.Cells.Item(21, 17).Value = 3.5
$ws.Cells.Item(21, 21).Value = 9.5
$ws.Cells.Item(21, 22).Value = 8.5
$ws.Cells.Item(21, 25).Value = 21
$ws.Cells.Item(21, 33).Value = 17

# Row 22
$ws.Cells.Item(22, 7).Value = 3.7
$ws.Cells.Item(22, 9).Value = 1.8
$ws.Cells.Item(22, 10).Value = 1.02
$ws.Cells.Item(22, 12).Value = 1.11
$ws.Cells.Item(22, 20).Value = 21
$ws.Cells.Item(22, 21).Value = 26
$ws.Cells.Item(22, 22).Value = 15
$ws.Cells.Item(22, 24).Value = 26
$ws.Cells.Item(22, 25).Value = 26
$ws.Cells.Item(22, 26).Value = 23
$ws.Cells.Item(22, 27).Value = 9

Write-Host "Applied 221 cell updates"